$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells that become empty (cell removed from sparse XML) - clear individually
$ws.Range("F5").ClearContents()
$ws.Range("G5").ClearContents()
$ws.Range("H11").ClearContents()
$ws.Range("I11").ClearContents()
$ws.Range("F32").ClearContents()
$ws.Range("G32").ClearContents()
$ws.Range("F38").ClearContents()
$ws.Range("G38").ClearContents()
$ws.Range("H53").ClearContents()
$ws.Range("I53").ClearContents()
$ws.Range("F104").ClearContents()
$ws.Range("G104").ClearContents()

# Set numeric cell values
$ws.Range("J4").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("J7").Value = 913
$ws.Range("J22").Value = 0
$ws.Range("G25").Value = 7261
$ws.Range("I25").Value = 1
$ws.Range("H32").Value = 0
$ws.Range("G33").Value = 8603
$ws.Range("I33").Value = 1
$ws.Range("K36").Value = 0
$ws.Range("H38").Value = 0
$ws.Range("J49").Value = 611
$ws.Range("K49").Value = 1410
$ws.Range("J52").Value = 1319
$ws.Range("K52").Value = 0
$ws.Range("G58").Value = 6941
$ws.Range("I58").Value = 2
$ws.Range("J64").Value = 0
$ws.Range("J68").Value = 812
$ws.Range("I69").Value = 1
$ws.Range("K71").Value = 665
$ws.Range("J72").Value = 1610
$ws.Range("J78").Value = 1748
$ws.Range("J79").Value = 1339
$ws.Range("H81").Value = 3
$ws.Range("J81").Value = 1906
$ws.Range("J83").Value = 319
$ws.Range("K83").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("J96").Value = 717
$ws.Range("J97").Value = 1889
$ws.Range("J100").Value = 747
$ws.Range("I104").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("K141").Value = 483
$ws.Range("J142").Value = 447
$ws.Range("J144").Value = 1019
$ws.Range("K144").Value = 0
$ws.Range("K145").Value = 0
$ws.Range("I146").Value = 2
$ws.Range("J151").Value = 0
$ws.Range("J152").Value = 660
$ws.Range("K155").Value = 0
$ws.Range("J170").Value = 0
$ws.Range("J171").Value = 0
$ws.Range("J172").Value = 0
$ws.Range("J173").Value = 503
$ws.Range("K175").Value = 2127
$ws.Range("J176").Value = 395
$ws.Range("J177").Value = 0
$ws.Range("K180").Value = 0
$ws.Range("J181").Value = 1747
$ws.Range("J184").Value = 0
$ws.Range("J187").Value = 0
$ws.Range("K189").Value = 0

# Set string cell values
$ws.Range("F25").Value = "incongruent"
$ws.Range("F33").Value = "incongruent"
$ws.Range("F58").Value = "incongruent"
